$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New arrival row (row 8) - Friday, Jan 13 LOT flight LO3507 from Krakow,
# same aircraft type as the earlier LOT flight but in Star Alliance livery.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Friday, Jan 13"
$ws.Range("C8").Value = "5:55 PM"
$ws.Range("D8").Value = "LO3507"
$ws.Range("E8").Value = "Krakow"
$ws.Range("F8").Value = "(KRK)"
$ws.Range("G8").Value = "LOT (Star Alliance Livery) "
$ws.Range("H8").Value = "E75S"
$ws.Range("I8").Value = "(SP-LIO)"
$ws.Range("J8").Value = "5:24 PM"
$ws.Range("L8").Value = "0 hours, -31 minutes"

# K8/M8 stay blank, but every other data row has an (empty) cell present in
# those two spacer columns, so mirror that by copying the existing blank
# formatting from the row above instead of leaving the cells absent.
$ws.Range("K7").Copy($ws.Range("K8"))
$ws.Range("M7").Copy($ws.Range("M8"))
